# Auto-generated Excel COM-interop script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to remain
# text (matching the original inlineStr cell type) by pre-formatting as text.
$textCells = @("D5", "D6", "D8", "D10", "D11", "D15", "D16", "D18", "D22", "D24", "D25", "D27", "D30", "D31", "D32", "D34", "D35", "D39", "D40", "D41", "D47", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "28.108.87"
$ws.Range("E2").Value = "  +2.41%  "
$ws.Range("D3").Value = "1.653.68"
$ws.Range("E3").Value = "  +2.35%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "214.15"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").Value = "0.530"
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "23.61"
$ws.Range("E8").Value = "  +3.68%  "
$ws.Range("E9").Value = "  +2.07%  "
$ws.Range("D10").Value = "0.0615"
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("D11").Value = "0.0875"
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("D12").Value = "1.888.88"
$ws.Range("E12").Value = "  +2.43%  "
$ws.Range("D13").Value = "1.654.35"
$ws.Range("E13").Value = "  +2.36%  "
$ws.Range("E14").Value = "  +1.47%  "
$ws.Range("D15").Value = "0.569"
$ws.Range("E15").Value = "  +3.65%  "
$ws.Range("D16").Value = "65.85"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("D17").Value = "28.102.19"
$ws.Range("E17").Value = "  +2.42%  "
$ws.Range("D18").Value = "233.74"
$ws.Range("E18").Value = "  +1.32%  "
$ws.Range("E19").Value = "  +2.45%  "
$ws.Range("D20").Value = "0.0₃0726"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").Value = "10.69"
$ws.Range("E22").Value = "  +5.46%  "
$ws.Range("D24").Value = "2.15"
$ws.Range("E24").Value = "  +3.58%  "
$ws.Range("D25").Value = "152.52"
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("D27").Value = "15.84"
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "1.20"
$ws.Range("E30").Value = "  +2.09%  "
$ws.Range("D31").Value = "0.0486"
$ws.Range("E31").Value = "  +0.69%  "
$ws.Range("D32").Value = "3.35"
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("D33").Value = "1.453.91"
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("D34").Value = "3.09"
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("D35").Value = "1.58"
$ws.Range("E35").Value = "  +3.05%  "
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("E37").Value = "  +4.33%  "
$ws.Range("E38").Value = "  +1.90%  "
$ws.Range("D39").Value = "0.933"
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("D40").Value = "0.562"
$ws.Range("E40").Value = "  +0.91%  "
$ws.Range("D41").Value = "69.48"
$ws.Range("E41").Value = "  +2.29%  "
$ws.Range("E42").Value = "  +3.78%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("E45").Value = "  +6.21%  "
$ws.Range("E46").Value = "  +1.50%  "
$ws.Range("D47").Value = "5.43"
$ws.Range("E47").Value = "  +3.48%  "
$ws.Range("D48").Value = "1.796.96"
$ws.Range("E48").Value = "  +2.19%  "
$ws.Range("D49").Value = "89.19"
$ws.Range("E49").Value = "  +3.15%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0105"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.102"
$ws.Range("E51").Value = "  +0.95%  "

# Restore default (Normal) style on the text-forced cells so no stray
# number-format styling is left behind on them.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
